$d = $word.ActiveDocument

# -----------------------------------------------------------------
# The paragraph "J'ai un pb avec mes tests -> Ils ne marchent pas,
# alors que chez Clérice ils marchent très bien [emoji]" (paragraph
# index 8, right after "Penser à faire les tests depuis une bd
# fermée, sinon ça plante") is replaced by two brand-new paragraphs:
#   1. "Mon 1er test fonctionnel -> Faudra maintenant améliorer tout
#      ça : les généraliser + vérifier qu'on parle bien de la même
#      chose (pk None pour info_personnelles ?)" -- with "er"
#      rendered as superscript.
#   2. "Penser à nettoyer un peu mon code de mes notes devenues
#      inutiles"
# Every paragraph that follows keeps its original content; it is
# simply pushed two slots further down the document.
# -----------------------------------------------------------------

# 1) Remove the old paragraph entirely (text + the emoji run) by
#    deleting its whole range, including the trailing paragraph mark.
$oldPara = $d.Paragraphs.Item(8)
$oldPara.Range.Delete() | Out-Null

# 2) Insert the two replacement paragraphs right after paragraph 7
#    ("Penser à faire les tests depuis une bd fermée, sinon ça
#    plante"), which is a plain paragraph (no list style) -- so the
#    new paragraphs inherit plain formatting too.
$anchor = $d.Paragraphs.Item(7)
$anchor.Range.InsertParagraphAfter() | Out-Null

$para1 = $d.Paragraphs.Item(8)
$para1.Range.Text = "Mon 1er test fonctionnel -> Faudra maintenant améliorer tout ça : les généraliser + vérifier qu’on parle bien de la même chose (pk None pour info_personnelles ?)"

# 3) Insert the second new paragraph right after the first one (while
#    the first paragraph's text is still plain, so the paragraph mark
#    used as the insertion anchor carries no stray superscript state).
$para1 = $d.Paragraphs.Item(8)
$para1.Range.InsertParagraphAfter() | Out-Null
$para2 = $d.Paragraphs.Item(9)
$para2.Range.Text = "Penser à nettoyer un peu mon code de mes notes devenues inutiles"

# 4) Now make the "er" right after "Mon 1" superscript.
$para1 = $d.Paragraphs.Item(8)
$p1Start = $para1.Range.Start
$supStart = $p1Start + 5
$supEnd = $supStart + 2
$supRange = $d.Range($supStart, $supEnd)
$supRange.Font.Superscript = $true

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
